# Generate Report for Handback
# Refresh the generation/handoff/handback timestamps on the handback status
# report, as would happen when the report generator re-runs a bit later.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the first row ---
# This text is also shared (same shared-string entry) with the "Correspond
# Handoff Datetime" cell on the de-de sheet (H2), so both must be updated
# together to keep them in sync.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-06 21:29:14"

# --- zh-cn sheet: handoff / handback datetimes for the first data row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-06 21:29:00"
$wsZhCn.Range("K2").Value = "2016-09-06 21:29:33"

# --- de-de sheet: handoff / handback datetimes for the first data row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-06 21:29:14"
$wsDeDe.Range("K2").Value = "2016-09-06 21:29:41"
